$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.074.13"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").Value = "1.799.84"
$ws.Range("E3").Value = "  -1.23%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.32"
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5098"
$ws.Range("E7").Value = "  -2.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3904"
$ws.Range("E8").Value = "  +1.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07765"
$ws.Range("E9").Value = "  -3.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.100"
$ws.Range("E10").Value = "  -1.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.00"
$ws.Range("E11").Value = "  -2.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.324"
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.005"
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.22"
$ws.Range("E14").Value = "  -3.30%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.800.43"
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.281"
$ws.Range("E16").Value = "  -2.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.13"
$ws.Range("E17").Value = "  -2.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001070"
$ws.Range("E18").Value = "  -3.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06574"
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.20"
$ws.Range("E21").Value = "  -2.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.962"
$ws.Range("E22").Value = "  -0.99%  "
$ws.Range("D23").Value = "28.146.95"
$ws.Range("E23").Value = "  -1.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.07"
$ws.Range("E24").Value = "  -2.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.232"
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.39"
$ws.Range("E26").Value = "  +1.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.428"
$ws.Range("E27").Value = "  +0.81%  "
$ws.Range("D28").Value = "2.011.59"
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.35"
$ws.Range("E29").Value = "  -2.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.99"
$ws.Range("E30").Value = "  +2.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1088"
$ws.Range("E31").Value = "  -2.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.049"
$ws.Range("E32").Value = "  -2.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.655"
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.523"
$ws.Range("E34").Value = "  -2.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07027"
$ws.Range("E35").Value = "  -3.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.084"
$ws.Range("E36").Value = "  +3.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02337"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2156"
$ws.Range("E38").Value = "  -2.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.47"
$ws.Range("E39").Value = "  -6.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.978"
$ws.Range("E40").Value = "  -2.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6118"
$ws.Range("E41").Value = "  -3.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.003"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.154"
$ws.Range("E43").Value = "  -2.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.21"
$ws.Range("E44").Value = "  -1.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.303"
$ws.Range("E45").Value = "  -5.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5897"
$ws.Range("E46").Value = "  -3.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.715"
$ws.Range("E47").Value = "  -2.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.57"
$ws.Range("E48").Value = "  -2.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.203"
$ws.Range("E49").Value = "  -1.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.909"
$ws.Range("E50").Value = "  -2.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06752"
$ws.Range("E51").Value = "  -2.11%  "
